$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1819.2667
$ws.Range("I113").Value = 1499.5
$ws.Range("J113").Value = 1868.4615
$ws.Range("K113").Value = 1499.5
$ws.Range("L113").Value = 1868.4615
$ws.Range("M113").Value = 1754.5
$ws.Range("N113").Value = -8376.461499999999
$ws.Range("H131").Value = 24163.592
$ws.Range("I131").Value = 29438.828
$ws.Range("J131").Value = 3648.7778
$ws.Range("K131").Value = 88316.484
$ws.Range("L131").Value = 10946.3334
$ws.Range("M131").Value = -83276.484
$ws.Range("N131").Value = -21026.3334
$ws.Range("H137").Value = 42220.24
$ws.Range("I137").Value = 1936.5454
$ws.Range("J137").Value = 73871.71000000001
$ws.Range("K137").Value = 5809.6362
$ws.Range("L137").Value = 221615.13
$ws.Range("M137").Value = -3259.6362
$ws.Range("N137").Value = -226715.13

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1804.72
$ws.Range("I2").Value = 1706.9474
$ws.Range("J2").Value = 2114.3333
$ws.Range("K2").Value = 1706.9474
$ws.Range("L2").Value = 2114.3333
$ws.Range("M2").Value = -1593.9474
$ws.Range("N2").Value = -2340.3333
$ws.Range("H61").Value = 3233.3333
$ws.Range("I61").Value = 1404
$ws.Range("J61").Value = 8721.333000000001
$ws.Range("K61").Value = 1404
$ws.Range("L61").Value = 8721.333000000001
$ws.Range("M61").Value = -1192
$ws.Range("N61").Value = -9145.333000000001
$ws.Range("H116").Value = 1804.72
$ws.Range("I116").Value = 1706.9474
$ws.Range("J116").Value = 2114.3333
$ws.Range("K116").Value = 1706.9474
$ws.Range("L116").Value = 2114.3333
$ws.Range("M116").Value = 587.0526
$ws.Range("N116").Value = -6702.3333
$ws.Range("H136").Value = 3233.3333
$ws.Range("I136").Value = 1404
$ws.Range("J136").Value = 8721.333000000001
$ws.Range("K136").Value = 4212
$ws.Range("L136").Value = 26163.999
$ws.Range("M136").Value = -1662
$ws.Range("N136").Value = -31263.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1804.72
$ws.Range("I3").Value = 1706.9474
$ws.Range("J3").Value = 2114.3333
$ws.Range("K3").Value = 1706.9474
$ws.Range("L3").Value = 2114.3333
$ws.Range("M3").Value = -1592.9474
$ws.Range("N3").Value = -2342.3333
$ws.Range("H105").Value = 2701.91
$ws.Range("I105").Value = 1339.2354
$ws.Range("J105").Value = 2981.012
$ws.Range("K105").Value = 1339.2354
$ws.Range("L105").Value = 2981.012
$ws.Range("M105").Value = 407.7646
$ws.Range("N105").Value = -6475.012000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 65.26667
$ws.Range("I7").Value = 46.166668
$ws.Range("J7").Value = 78
$ws.Range("K7").Value = 46.166668
$ws.Range("L7").Value = 78
$ws.Range("M7").Value = 66.833332
$ws.Range("N7").Value = -304
$ws.Range("H31").Value = 3630.1
$ws.Range("I31").Value = 2481.1667
$ws.Range("J31").Value = 5353.5
$ws.Range("K31").Value = 2481.1667
$ws.Range("L31").Value = 5353.5
$ws.Range("M31").Value = -2186.1667
$ws.Range("N31").Value = -5943.5
$ws.Range("H34").Value = 3630.1
$ws.Range("I34").Value = 2481.1667
$ws.Range("J34").Value = 5353.5
$ws.Range("K34").Value = 2481.1667
$ws.Range("L34").Value = 5353.5
$ws.Range("M34").Value = -2279.1667
$ws.Range("N34").Value = -5757.5
$ws.Range("H107").Value = 683.9697
$ws.Range("I107").Value = 496.3684
$ws.Range("J107").Value = 938.5714
$ws.Range("K107").Value = 496.3684
$ws.Range("L107").Value = 938.5714
$ws.Range("M107").Value = 1423.6316
$ws.Range("N107").Value = -4778.5714
$ws.Range("H110").Value = 61580.4
$ws.Range("J110").Value = 61580.4
$ws.Range("L110").Value = 61580.4
$ws.Range("N110").Value = -69760.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 239.375
$ws.Range("I23").Value = 217.6
$ws.Range("K23").Value = 652.8
$ws.Range("M23").Value = -417.8
$ws.Range("H131").Value = 809.6486
$ws.Range("J131").Value = 1002.8095
$ws.Range("L131").Value = 3008.4285
$ws.Range("N131").Value = -13088.4285
$ws.Range("H132").Value = 674680
$ws.Range("I132").Value = 1427.1428
$ws.Range("J132").Value = 1263776.2
$ws.Range("K132").Value = 12844.2852
$ws.Range("L132").Value = 11373985.8
$ws.Range("M132").Value = -10314.2852
$ws.Range("N132").Value = -11379045.8
$ws.Range("H134").Value = 3939.6333
$ws.Range("I134").Value = 1368.3846
$ws.Range("J134").Value = 5905.8823
$ws.Range("K134").Value = 4105.1538
$ws.Range("L134").Value = 17717.6469
$ws.Range("M134").Value = 964.8462
$ws.Range("N134").Value = -27857.6469
$ws.Range("H141").Value = 5886540.5
$ws.Range("I141").Value = 11112799
$ws.Range("J141").Value = 7000
$ws.Range("K141").Value = 33338397
$ws.Range("L141").Value = 21000
$ws.Range("M141").Value = -33333217
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 134.24138
$ws.Range("I107").Value = 139.20833
$ws.Range("K107").Value = 139.20833
$ws.Range("M107").Value = 1780.79167

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 260.31033
$ws.Range("I55").Value = 278.25
$ws.Range("J55").Value = 220.44444
$ws.Range("K55").Value = 278.25
$ws.Range("L55").Value = 220.44444
$ws.Range("M55").Value = -105.25
$ws.Range("N55").Value = -566.44444
$ws.Range("H61").Value = 8348.385
$ws.Range("I61").Value = 20385.8
$ws.Range("J61").Value = 825
$ws.Range("K61").Value = 20385.8
$ws.Range("L61").Value = 825
$ws.Range("M61").Value = -20183.8
$ws.Range("N61").Value = -1229
$ws.Range("H113").Value = 8348.385
$ws.Range("I113").Value = 20385.8
$ws.Range("J113").Value = 825
$ws.Range("K113").Value = 20385.8
$ws.Range("L113").Value = 825
$ws.Range("M113").Value = -18215.8
$ws.Range("N113").Value = -5165

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 684.9286
$ws.Range("I107").Value = 730.75
$ws.Range("J107").Value = 623.8333
$ws.Range("K107").Value = 2192.25
$ws.Range("L107").Value = 1871.4999
$ws.Range("M107").Value = -272.25
$ws.Range("N107").Value = -5711.4999
$ws.Range("H113").Value = 317.0625
$ws.Range("I113").Value = 300.18182
$ws.Range("J113").Value = 354.2
$ws.Range("K113").Value = 900.54546
$ws.Range("L113").Value = 1062.6
$ws.Range("M113").Value = 1269.45454
$ws.Range("N113").Value = -5402.6
